# Alarm Normal load method changes
# Add two new columns (L, M) with header + data cells to the "Add Panels" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Panels")

# --- New header cells (row 7) -------------------------------------------------
$ws.Range("L7").Value = "AlarmLoadingDetail"
$ws.Range("K7").Copy()
$ws.Range("L7").PasteSpecial(-4122)   # xlPasteFormats - copy K7's header style

$ws.Range("M7").Value = "StandbyLoadingDetail"
$ws.Range("K7").Copy()
$ws.Range("M7").PasteSpecial(-4122)   # xlPasteFormats

# --- New data cells (row 8) ---------------------------------------------------
$ws.Range("L8").Value = "Battery Alarm (A)"
$ws.Range("B8").Copy()
$ws.Range("L8").PasteSpecial(-4122)   # xlPasteFormats - copy B8's data style

$ws.Range("M8").Value = "Battery Standby (A)"
$ws.Range("B8").Copy()
$ws.Range("M8").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = $false

# --- Column widths for the two new columns ------------------------------------
$ws.Columns.Item(12).ColumnWidth = 16.833333333333336
$ws.Columns.Item(13).ColumnWidth = 18.833333333333336

# --- View state: scroll and selection -----------------------------------------
$win = $wb.Windows.Item(1)
$win.ScrollColumn = 3
$ws.Range("L15").Select()
